# Time Tables.xlsx update
# Adds 6 new Pygame-related log rows (12-17), blanks out / extends the
# trailing grid down to row 30, and adds a new "Full Total" labeled
# grand-total row (29/30), replacing the old F20 sub-total.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 0) Grab the format for the new grand-total formula cell (F30) from the
#    old F20 sub-total cell before we touch/clear it.
# ---------------------------------------------------------------------
$ws.Range("F20").Copy() | Out-Null
$ws.Range("F30").PasteSpecial(-4122) | Out-Null       # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 1) New data rows 12-16 : Title / Date / Time Started / Time Ended / Total
# ---------------------------------------------------------------------

# -- Title column (A12:A16) : same look as the existing title cells (A2) --
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A12:A16").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# -- Date column (B12:B17) : same look as existing date cells, but the
#    new rows use a plain (non right-aligned) date layout --
$ws.Range("C4").Copy() | Out-Null
$ws.Range("B12:B17").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("B12:B17").NumberFormat = "m/d/yy"

# -- Time Started / Time Ended columns (C12:D17) --
$ws.Range("C4:D4").Copy() | Out-Null
$ws.Range("C12:D17").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# -- Total column (E12:E17) --
$ws.Range("E4").Copy() | Out-Null
$ws.Range("E12:E17").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

# -- values --
$ws.Range("A12").Value = "Trying to Redo everything in Pygame"
$ws.Range("B12").Value = 44319
$ws.Range("C12").Value = 0.49722222222222223
$ws.Range("D12").Value = 0.5263888888888889
$ws.Range("E12").Value = 0.02900462962962963

$ws.Range("A13").Value = "Sprites in Pygame"
$ws.Range("B13").Value = 44319
$ws.Range("C13").Value = 0.9270833333333334
$ws.Range("D13").Value = 0.08055555555555556
$ws.Range("E13").Value = 0.15314814814814814

$ws.Range("A14").Value = "Loading Door Game Concept Test"
$ws.Range("B14").Value = 44320
$ws.Range("C14").Value = 0.6298611111111111
$ws.Range("D14").Value = 0.7458333333333333
$ws.Range("E14").Value = 0.11597222222222223

$ws.Range("A15").Value = "Pygame GUI"
$ws.Range("B15").Value = 44321
$ws.Range("C15").Value = 0.6652777777777777
$ws.Range("D15").Value = 0.8201388888888889
$ws.Range("E15").Value = 0.15493055555555554

$ws.Range("A16").Value = "Collisions???"
$ws.Range("B16").Value = 44321
$ws.Range("C16").Value = 0.8208333333333333
$ws.Range("D16").Value = 0.9715277777777778
$ws.Range("E16").Value = 0.15083333333333335

# ---------------------------------------------------------------------
# 2) Row 17 : "Redoing everything in Pygame" - special white, left
#    aligned, borderless title cell, rest of the row like rows 12-16
# ---------------------------------------------------------------------
$r = $ws.Range("A17")
$r.Interior.Color = 16777215        # white fill (matches the other white cells)
$r.Font.Name = "Arial"
$r.Font.Color = 0                   # explicit black
$r.Borders.LineStyle = -4142        # xlLineStyleNone - remove the box border
$r.HorizontalAlignment = -4131      # xlLeft
$r.Value = "Redoing everything in Pygame"

$ws.Range("B17").Value = 44322
$ws.Range("C17").Value = 0.059722222222222225
$ws.Range("D17").Value = 0.1326388888888889
$ws.Range("E17").Value = 0.07291666666666667

# ---------------------------------------------------------------------
# 3) Remove the old sub-total formula that used to live in F20
# ---------------------------------------------------------------------
$ws.Range("F20").Clear() | Out-Null

# ---------------------------------------------------------------------
# 4) Extend the blank bordered grid through row 30 (rows 18-19 already
#    carry the right blank look; rows 20-28 get the same treatment)
# ---------------------------------------------------------------------
$ws.Range("A18:E19").Copy() | Out-Null
$ws.Range("A20:E28").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Row 20 fully matches the plain "D-column" blank look (no special total col)
$ws.Range("D18").Copy() | Out-Null
$ws.Range("A20:E20").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Row 21 E column also follows the plain look (old bottom-border-only E21 cell)
$ws.Range("D18").Copy() | Out-Null
$ws.Range("E21").PasteSpecial(-4122) | Out-Null       # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A18:E19").Copy() | Out-Null
$ws.Range("A29:E30").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 5) New grand total row: "Full Total" label (green, bold) + formula
# ---------------------------------------------------------------------
$ws.Range("A1").Copy() | Out-Null
$ws.Range("F29").PasteSpecial(-4122) | Out-Null       # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("F29").Font.Color = 0                       # explicit black, bold stays
$ws.Range("F29").Value = "Full Total"

$ws.Range("F30").Formula = "=SUM(E2:E30)"
